# Scheduled-runner market data refresh for Brynhildr_Profits
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N)
# for the leves whose backing market prices moved since the last run.

$wb = $excel.ActiveWorkbook

# ================= Sheet: ALC =================
$ws = $wb.Worksheets.Item("ALC")

# Row 12
$ws.Range("H12").Value = 107.09091
$ws.Range("I12").Value = 104.22222
$ws.Range("K12").Value = 104.22222
$ws.Range("M12").Value = 65.77778000000001
# Row 62
$ws.Range("H62").Value = 2482.2
$ws.Range("I62").Value = 2135
$ws.Range("K62").Value = 2135
$ws.Range("M62").Value = -1511
# Row 65
$ws.Range("H65").Value = 2482.2
$ws.Range("I65").Value = 2135
$ws.Range("K65").Value = 10675
$ws.Range("M65").Value = -7555
# Row 99
$ws.Range("H99").Value = 326.66666
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
# Row 118
$ws.Range("H118").Value = 450
$ws.Range("I118").Value = 450
$ws.Range("K118").Value = 1350
$ws.Range("M118").Value = 307
# Row 119
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("M119").ClearContents()
# Row 132
$ws.Range("H132").Value = 6427.2666
$ws.Range("I132").Value = 7371.879
$ws.Range("K132").Value = 22115.637
$ws.Range("M132").Value = -19585.637
# Row 138
$ws.Range("H138").Value = 2111.4119
$ws.Range("I138").Value = 982.05554
$ws.Range("K138").Value = 2946.16662
$ws.Range("M138").Value = 2193.83338

# ================= Sheet: ARM =================
$ws = $wb.Worksheets.Item("ARM")

# Row 28
$ws.Range("H28").Value = 14049.333
$ws.Range("I28").Value = 7194.25
$ws.Range("J28").Value = 27759.5
$ws.Range("K28").Value = 7194.25
$ws.Range("L28").Value = 27759.5
$ws.Range("M28").Value = -7002.25
$ws.Range("N28").Value = -28143.5
# Row 63
$ws.Range("H63").Value = 2191.1667
$ws.Range("I63").Value = 2202.5881
$ws.Range("K63").Value = 2202.5881
$ws.Range("M63").Value = -1516.5881
# Row 66
$ws.Range("H66").Value = 2191.1667
$ws.Range("I66").Value = 2202.5881
$ws.Range("K66").Value = 11012.9405
$ws.Range("M66").Value = -7580.940500000001
# Row 99
$ws.Range("H99").Value = 14049.333
$ws.Range("I99").Value = 7194.25
$ws.Range("J99").Value = 27759.5
$ws.Range("K99").Value = 7194.25
$ws.Range("L99").Value = 27759.5
$ws.Range("M99").Value = -4199.25
$ws.Range("N99").Value = -33749.5
# Row 122
$ws.Range("H122").Value = 3127.7144
$ws.Range("I122").Value = 2973.75
$ws.Range("J122").Value = 3333
$ws.Range("K122").Value = 8921.25
$ws.Range("L122").Value = 9999
$ws.Range("M122").Value = -6471.25
$ws.Range("N122").Value = -14899

# ================= Sheet: BSM =================
$ws = $wb.Worksheets.Item("BSM")

# Row 99
$ws.Range("H99").Value = 5993.952
$ws.Range("I99").Value = 7101.1177
$ws.Range("J99").Value = 1288.5
$ws.Range("K99").Value = 7101.1177
$ws.Range("L99").Value = 1288.5
$ws.Range("M99").Value = -5603.1177
$ws.Range("N99").Value = -4284.5
# Row 105
$ws.Range("H105").Value = 2210.653
$ws.Range("I105").Value = 1532.2106
$ws.Range("K105").Value = 1532.2106
$ws.Range("M105").Value = 214.7893999999999
# Row 112
$ws.Range("H112").Value = 50000
$ws.Range("J112").Value = 50000
$ws.Range("L112").Value = 50000
$ws.Range("N112").Value = -52954

# ================= Sheet: CRP =================
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 2475.25
$ws.Range("I31").Value = 2405.8635
$ws.Range("J31").Value = 2627.9
$ws.Range("K31").Value = 2405.8635
$ws.Range("L31").Value = 2627.9
$ws.Range("M31").Value = -2110.8635
$ws.Range("N31").Value = -3217.9
# Row 34
$ws.Range("H34").Value = 2475.25
$ws.Range("I34").Value = 2405.8635
$ws.Range("J34").Value = 2627.9
$ws.Range("K34").Value = 2405.8635
$ws.Range("L34").Value = 2627.9
$ws.Range("M34").Value = -2203.8635
$ws.Range("N34").Value = -3031.9
# Row 119
$ws.Range("H119").Value = 90000
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 90000
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 90000
$ws.Range("M119").ClearContents()
$ws.Range("N119").Value = -99676
# Row 122
$ws.Range("H122").Value = 1773.6216
$ws.Range("I122").Value = 1475.7693
$ws.Range("J122").Value = 2477.6365
$ws.Range("K122").Value = 4427.3079
$ws.Range("L122").Value = 7432.9095
$ws.Range("M122").Value = -1977.3079
$ws.Range("N122").Value = -12332.9095
# Row 134
$ws.Range("H134").Value = 2575
$ws.Range("I134").Value = 2432.5
$ws.Range("K134").Value = 7297.5
$ws.Range("M134").Value = -4762.5

# ================= Sheet: CUL =================
$ws = $wb.Worksheets.Item("CUL")

# Row 131
$ws.Range("H131").Value = 3185.6
$ws.Range("J131").Value = 3185.6
$ws.Range("L131").Value = 9556.799999999999
$ws.Range("N131").Value = -19636.8
# Row 132
$ws.Range("H132").Value = 788.6429000000001
$ws.Range("I132").Value = 574.25
$ws.Range("J132").Value = 1074.5
$ws.Range("K132").Value = 5168.25
$ws.Range("L132").Value = 9670.5
$ws.Range("M132").Value = -2638.25
$ws.Range("N132").Value = -14730.5

# ================= Sheet: GSM =================
$ws = $wb.Worksheets.Item("GSM")

# Row 80
$ws.Range("H80").Value = 1679.25
$ws.Range("I80").Value = 1523.5
$ws.Range("J80").Value = 1835
$ws.Range("K80").Value = 1523.5
$ws.Range("L80").Value = 1835
$ws.Range("M80").Value = -525.5
$ws.Range("N80").Value = -3831
# Row 83
$ws.Range("H83").Value = 1679.25
$ws.Range("I83").Value = 1523.5
$ws.Range("J83").Value = 1835
$ws.Range("K83").Value = 7617.5
$ws.Range("L83").Value = 9175
$ws.Range("M83").Value = -2625.5
$ws.Range("N83").Value = -19159
# Row 122
$ws.Range("H122").Value = 66803.875
$ws.Range("I122").Value = 94907.09
$ws.Range("K122").Value = 284721.27
$ws.Range("M122").Value = -282271.27
# Row 132
$ws.Range("H132").Value = 10421.873
$ws.Range("I132").Value = 12211.048
$ws.Range("K132").Value = 36633.144
$ws.Range("M132").Value = -34103.144
# Row 136
$ws.Range("H136").Value = 75954.664
$ws.Range("J136").Value = 75954.664
$ws.Range("L136").Value = 227863.992
$ws.Range("N136").Value = -232963.992

# ================= Sheet: LTW =================
$ws = $wb.Worksheets.Item("LTW")

# Row 25
$ws.Range("H25").Value = 1600
$ws.Range("J25").Value = 1600
$ws.Range("L25").Value = 1600
$ws.Range("N25").Value = -2060
# Row 40
$ws.Range("H40").Value = 3500.5
$ws.Range("I40").Value = 2002.6666
$ws.Range("K40").Value = 2002.6666
$ws.Range("M40").Value = -1866.6666
# Row 122
$ws.Range("H122").Value = 6800.1665
$ws.Range("I122").Value = 5099.5
$ws.Range("J122").Value = 8500.833000000001
$ws.Range("K122").Value = 15298.5
$ws.Range("L122").Value = 25502.499
$ws.Range("M122").Value = -12848.5
$ws.Range("N122").Value = -30402.499

# ================= Sheet: WVR =================
$ws = $wb.Worksheets.Item("WVR")

# Row 4
$ws.Range("H4").Value = 92
$ws.Range("J4").Value = 100
$ws.Range("L4").Value = 100
$ws.Range("N4").Value = -326
# Row 107
$ws.Range("H107").Value = 1815.3103
$ws.Range("I107").Value = 1147.7273
$ws.Range("J107").Value = 3913.4285
$ws.Range("K107").Value = 3443.1819
$ws.Range("L107").Value = 11740.2855
$ws.Range("M107").Value = -1523.1819
$ws.Range("N107").Value = -15580.2855

